$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.924.96'
$ws.Cells.Item(2, 5).Value = '  +0.04%  '

$ws.Cells.Item(3, 4).Value = '3.416.93'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '

$ws.Cells.Item(4, 5).Value = '  +0.34%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '409.62'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.30%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '128.23'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.48%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.633'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +7.28%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.05%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.730'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +6.39%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +10.75%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '42.50'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +1.55%  '

$ws.Cells.Item(12, 5).Value = '  +0.14%  '

$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '3.959.18'
$ws.Cells.Item(13, 5).Value = '  +0.11%  '

$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '9.02'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +7.42%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '21.17'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +7.15%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '0.0000203'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +45.72%  '

$ws.Cells.Item(17, 4).Value = '3.404.90'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '

$ws.Cells.Item(18, 5).Value = '  +5.27%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '1.08'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +6.50%  '

$ws.Cells.Item(20, 4).Value = '61.938.35'
$ws.Cells.Item(20, 5).Value = '  +0.09%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '445.14'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +42.85%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '92.19'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +10.77%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '3.16'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.05%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '12.88'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +1.25%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '3.24'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +2.88%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '33.17'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +12.15%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '8.80'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +9.30%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '4.80'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.42%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.74'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.52%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '7.55'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.74%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '11.93'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +5.17%  '

$ws.Cells.Item(32, 5).Value = '  -1.38%  '

$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '0.113'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -1.32%  '

$ws.Cells.Item(34, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '42.66'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.44%  '

$ws.Cells.Item(35, 5).Value = '  -0.01%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.0495'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +2.83%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '53.40'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +3.96%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +0.11%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '3.37'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.10%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '0.133'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +7.29%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '2.91'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -2.03%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.312'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.97%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '140.95'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.81%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '4.21'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +6.86%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '1.97'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.42%  '

$ws.Cells.Item(46, 5).Value = '  +8.52%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '16.46'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.77%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '22.33'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +5.13%  '

$ws.Cells.Item(49, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(49, 4).Value = '3.767.92'
$ws.Cells.Item(49, 5).Value = '  +0.34%  '

$ws.Cells.Item(50, 2).Value = 'ThetaToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '2.08'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +6.86%  '

$ws.Cells.Item(51, 4).Value = '2.112.79'
$ws.Cells.Item(51, 5).Value = '  +0.30%  '
